$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update Correspond Handoff / Handback datetimes for the
# first file row (3b87fe38-...) to reflect the newly generated handback report.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 14:39:36"
$wsZhCn.Range("H2").Value = "2016-03-19 14:39:55"

# "de-de" sheet: same update for its first file row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 14:39:39"
$wsDeDe.Range("H2").Value = "2016-03-19 14:40:00"
